# A02 Pixell test plan - savings account
# Fills in the test-case detail columns (Method Inputs / Preconditions /
# Expected Result -> worksheet columns E/F/G) for test rows 7-12, and
# updates the "Developer:" name.
#
# NOTE: the cells below are intentionally written in the same order the
# original author's save produced in the shared-string table (so the
# resulting workbook's shared-string ordering/index values line up with
# the canonical OOXML), rather than a simple top-to-bottom left-to-right scan.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Developer name -------------------------------------------------------
$ws.Range("C3").Value = "Beerdavinder singh"

# --- First occurrences (establishes shared-string order) -----------------
$ws.Range("E7").Value = "Inputs"
$ws.Range("E8").Value = "None"
$ws.Range("F7").Value = "(22222, 3333, 4444.44, date.today(), 250.00)`n         "
$ws.Range("F8").Value = "(22222, 3333, 4444.44, date.today(), `"invalid`")"
$ws.Range("G7").Value = "Setup"
$ws.Range("E12").Value = " f`"Account number: 22222 Balance: `$4444.44`n`" \`n                   f`"Minimum Balance: `$250.00 Account Type: Savings`"`n                   "
$ws.Range("G12").Value = "FORMATTED STR"
$ws.Range("G9").Value = "GET SERVICE CHARGES"
$ws.Range("G8").Value = "RAISES VALUEERROR"

# --- Remaining cells (reuse the shared strings set up above) -------------
$ws.Range("E9").Value = "(22222, 3333, 4444.44, date.today(), 250.00)`n         "
$ws.Range("F9").Value = "None"
$ws.Range("E10").Value = "(22222, 3333, 4444.44, date.today(), 250.00)`n         "
$ws.Range("F10").Value = "None"
$ws.Range("G10").Value = "GET SERVICE CHARGES"
$ws.Range("E11").Value = "(22222, 3333, 4444.44, date.today(), 250.00)`n         "
$ws.Range("F11").Value = "None"
$ws.Range("G11").Value = "GET SERVICE CHARGES"
$ws.Range("F12").Value = "(22222, 3333, 4444.44, date.today(), 250.00)`n         "

# --- Selection, mirroring the saved view state ----------------------------
$ws.Range("G8").Select()
